$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append a new row of measurements (row 7) below the existing data (rows 2-6)
# ---------------------------------------------------------------------------

# Date for column A (4/2/2025 -> Excel serial 45749)
$ws.Range("A7").Value2 = 45749

# Measurement values for columns B..M
$ws.Range("B7").Value2 = 35
$ws.Range("C7").Value2 = 40
$ws.Range("D7").Value2 = 38.1
$ws.Range("E7").Value2 = 37.9
$ws.Range("F7").Value2 = 39.200000000000003
$ws.Range("G7").Value2 = 32.4
$ws.Range("H7").Value2 = 31.3
$ws.Range("I7").Value2 = 38.200000000000003
$ws.Range("J7").Value2 = 36.299999999999997
$ws.Range("K7").Value2 = 39.299999999999997
$ws.Range("L7").Value2 = 31.7
$ws.Range("M7").Value2 = 39.200000000000003

# ---------------------------------------------------------------------------
# Formatting: mirror the formatting already used on the data block.
# Column A keeps the date format used by the rest of column A (copy from A2,
# the row that - like this new row - sits right under a border, so it only
# carries a bottom/left/right border, no top).
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Columns B:M get the same "no top border" box used on row 2 (the row
# immediately below the header), since row 7 sits directly below row 6's
# bottom border.
$ws.Range("B2:M2").Copy() | Out-Null
$ws.Range("B7:M7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Explicitly re-assert "no fill" on the new row so the format is fully
# self-contained (matches turning the fill off via Format Cells).
$ws.Range("B7:M7").Interior.ColorIndex = -4142   # xlColorIndexNone

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Selection state: the block B6:M7 is selected with M7 as the active cell.
# ---------------------------------------------------------------------------
$ws.Range("M7").Select() | Out-Null
$ws.Range("B6:M7").Select() | Out-Null
